$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Drop the 2020/2021 "Schliesstage" rows (rows 3 then 2, so indices stay valid)
$t.Rows.Item(3).Delete()
$t.Rows.Item(2).Delete()

# Rebuild the remaining row: first cell becomes the multi-line holiday list,
# second cell becomes the updated "Schliesstage" year label.
$cell1 = $t.Cell(1,1)
$cell1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Ostern</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>19.04. – 22.04.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:t>Brückentag</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>27.05.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>Pfingsten</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>13.06. – 17.06.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:t>Sommer</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>01.08. – 19.08.</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:t>Weihnachten</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>7</w:t></w:r><w:r><w:t>.12. – 3</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>.12.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cell1.Range.InsertXML($cell1Xml)

$cell2 = $t.Cell(1,2)
$cell2.Range.Text = "Schließtage 2022"
